$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H31").Value = 0
$ws.Range("I31").Value = 0
$ws.Range("K31").Value = 0
$ws.Range("M31").ClearContents()
$ws.Range("H74").Value = 71440140
$ws.Range("I74").Value = 125005120
$ws.Range("K74").Value = 125005120
$ws.Range("M74").Value = -125004184
$ws.Range("H77").Value = 71440140
$ws.Range("I77").Value = 125005120
$ws.Range("K77").Value = 625025600
$ws.Range("M77").Value = -625020920
$ws.Range("H103").Value = 1352.9375
$ws.Range("J103").Value = 1489.4615
$ws.Range("L103").Value = 4468.3845
$ws.Range("N103").Value = -5640.3845
$ws.Range("H132").Value = 1600
$ws.Range("I132").Value = 1639.2858
$ws.Range("J132").Value = 500
$ws.Range("K132").Value = 4917.857400000001
$ws.Range("L132").Value = 1500
$ws.Range("M132").Value = -2387.857400000001
$ws.Range("N132").Value = -6560
$ws.Range("H137").Value = 2040.6487
$ws.Range("I137").Value = 1928
$ws.Range("J137").Value = 2970
$ws.Range("K137").Value = 5784
$ws.Range("L137").Value = 8910
$ws.Range("M137").Value = -3234
$ws.Range("N137").Value = -14010
$ws.Range("H141").Value = 7248127.5
$ws.Range("I141").Value = 7938044.5
$ws.Range("K141").Value = 23814133.5
$ws.Range("M141").Value = -23808953.5

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1377196.1
$ws.Range("I32").Value = 1440050
$ws.Range("K32").Value = 1440050
$ws.Range("M32").Value = -1439763
$ws.Range("H61").Value = 4759.1187
$ws.Range("I61").Value = 2582.7659
$ws.Range("K61").Value = 2582.7659
$ws.Range("M61").Value = -2370.7659
$ws.Range("H74").Value = 18251.162
$ws.Range("I74").Value = 23491.482
$ws.Range("J74").Value = 4102.3
$ws.Range("K74").Value = 23491.482
$ws.Range("L74").Value = 4102.3
$ws.Range("M74").Value = -22617.482
$ws.Range("N74").Value = -5850.3
$ws.Range("H77").Value = 18251.162
$ws.Range("I77").Value = 23491.482
$ws.Range("J77").Value = 4102.3
$ws.Range("K77").Value = 117457.41
$ws.Range("L77").Value = 20511.5
$ws.Range("M77").Value = -113089.41
$ws.Range("N77").Value = -29247.5
$ws.Range("H88").Value = 2194.111
$ws.Range("I88").Value = 1948.75
$ws.Range("J88").Value = 2390.4
$ws.Range("K88").Value = 1948.75
$ws.Range("L88").Value = 2390.4
$ws.Range("M88").Value = -1542.75
$ws.Range("N88").Value = -3202.4
$ws.Range("H91").Value = 2194.111
$ws.Range("I91").Value = 1948.75
$ws.Range("J91").Value = 2390.4
$ws.Range("K91").Value = 1948.75
$ws.Range("L91").Value = 2390.4
$ws.Range("M91").Value = -544.75
$ws.Range("N91").Value = -5198.4
$ws.Range("H132").Value = 5087.1885
$ws.Range("I132").Value = 3554.4707
$ws.Range("K132").Value = 10663.4121
$ws.Range("M132").Value = -8133.4121
$ws.Range("H135").Value = 48597.5
$ws.Range("J135").Value = 48597.5
$ws.Range("L135").Value = 48597.5
$ws.Range("N135").Value = -58737.5
$ws.Range("H136").Value = 4759.1187
$ws.Range("I136").Value = 2582.7659
$ws.Range("K136").Value = 7748.297699999999
$ws.Range("M136").Value = -5198.297699999999

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 6174355.5
$ws.Range("I20").Value = 13890444
$ws.Range("J20").Value = 1484.6666
$ws.Range("K20").Value = 13890444
$ws.Range("L20").Value = 1484.6666
$ws.Range("M20").Value = -13890197
$ws.Range("N20").Value = -1978.6666
$ws.Range("H86").Value = 47671290
$ws.Range("I86").Value = 80752.766
$ws.Range("K86").Value = 80752.766
$ws.Range("M86").Value = -79629.766
$ws.Range("H89").Value = 47671290
$ws.Range("I89").Value = 80752.766
$ws.Range("K89").Value = 403763.83
$ws.Range("M89").Value = -398147.83
$ws.Range("H94").Value = 2568.9312
$ws.Range("I94").Value = 576.2
$ws.Range("J94").Value = 6997.222
$ws.Range("K94").Value = 576.2
$ws.Range("L94").Value = 6997.222
$ws.Range("M94").Value = -125.2
$ws.Range("N94").Value = -7899.222
$ws.Range("H107").Value = 40181976
$ws.Range("J107").Value = 4118.7
$ws.Range("L107").Value = 4118.7
$ws.Range("N107").Value = -7958.7
$ws.Range("H134").Value = 5137.5
$ws.Range("I134").Value = 1519.6
$ws.Range("J134").Value = 8755.4
$ws.Range("K134").Value = 4558.799999999999
$ws.Range("L134").Value = 26266.2
$ws.Range("M134").Value = -2023.799999999999
$ws.Range("N134").Value = -31336.2

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 6391.618
$ws.Range("J31").Value = 11097.6
$ws.Range("L31").Value = 11097.6
$ws.Range("N31").Value = -11687.6
$ws.Range("H34").Value = 6391.618
$ws.Range("J34").Value = 11097.6
$ws.Range("L34").Value = 11097.6
$ws.Range("N34").Value = -11501.6
$ws.Range("H58").Value = 8933086
$ws.Range("I58").Value = 13514854
$ws.Range("K58").Value = 13514854
$ws.Range("M58").Value = -13514651
$ws.Range("H82").Value = 32500
$ws.Range("J82").Value = 32500
$ws.Range("L82").Value = 32500
$ws.Range("N82").Value = -33222
$ws.Range("H85").Value = 32500
$ws.Range("J85").Value = 32500
$ws.Range("L85").Value = 32500
$ws.Range("N85").Value = -34996
$ws.Range("H99").Value = 9875.166999999999
$ws.Range("I99").Value = 11957.6
$ws.Range("K99").Value = 11957.6
$ws.Range("M99").Value = -10459.6
$ws.Range("H126").Value = 9875.166999999999
$ws.Range("I126").Value = 11957.6
$ws.Range("K126").Value = 35872.8
$ws.Range("M126").Value = -33402.8
$ws.Range("H132").Value = 4763.061
$ws.Range("I132").Value = 2062.2258
$ws.Range("J132").Value = 9414.5
$ws.Range("K132").Value = 6186.6774
$ws.Range("L132").Value = 28243.5
$ws.Range("M132").Value = -3656.6774
$ws.Range("N132").Value = -33303.5
$ws.Range("H134").Value = 6300.8237
$ws.Range("I134").Value = 1489.0625
$ws.Range("K134").Value = 4467.1875
$ws.Range("M134").Value = -1932.1875
$ws.Range("H136").Value = 8933086
$ws.Range("I136").Value = 13514854
$ws.Range("K136").Value = 40544562
$ws.Range("M136").Value = -40542012

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 210441100
$ws.Range("J4").Value = 200
$ws.Range("L4").Value = 600
$ws.Range("N4").Value = -824
$ws.Range("H38").Value = 104.55556
$ws.Range("I38").Value = 106.14286
$ws.Range("J38").Value = 99
$ws.Range("K38").Value = 318.42858
$ws.Range("L38").Value = 297
$ws.Range("M38").Value = 28.57141999999999
$ws.Range("N38").Value = -991
$ws.Range("H129").Value = 1731.8572
$ws.Range("I129").Value = 1010
$ws.Range("J129").Value = 2273.25
$ws.Range("K129").Value = 3030
$ws.Range("L129").Value = 6819.75
$ws.Range("M129").Value = 1970
$ws.Range("N129").Value = -16819.75

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 10479.571
$ws.Range("I70").Value = 8959.286
$ws.Range("J70").Value = 11999.857
$ws.Range("K70").Value = 8959.286
$ws.Range("L70").Value = 11999.857
$ws.Range("M70").Value = -8689.286
$ws.Range("N70").Value = -12539.857
$ws.Range("H73").Value = 10479.571
$ws.Range("I73").Value = 8959.286
$ws.Range("J73").Value = 11999.857
$ws.Range("K73").Value = 8959.286
$ws.Range("L73").Value = 11999.857
$ws.Range("M73").Value = -8023.286
$ws.Range("N73").Value = -13871.857
$ws.Range("H93").Value = 47951
$ws.Range("J93").Value = 47951
$ws.Range("L93").Value = 47951
$ws.Range("N93").Value = -51695
$ws.Range("H97").Value = 1044.5
$ws.Range("I97").Value = 957.80646
$ws.Range("K97").Value = 957.80646
$ws.Range("M97").Value = -461.80646
$ws.Range("H126").Value = 3143.2104
$ws.Range("I126").Value = 2501.7693
$ws.Range("K126").Value = 7505.3079
$ws.Range("M126").Value = -5035.3079
$ws.Range("H132").Value = 4766.5757
$ws.Range("I132").Value = 1700.4117
$ws.Range("K132").Value = 5101.2351
$ws.Range("M132").Value = -2571.2351

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H10").Value = 105
$ws.Range("I10").Value = 105
$ws.Range("J10").Value = 0
$ws.Range("K10").Value = 105
$ws.Range("L10").Value = 0
$ws.Range("M10").Value = 35
$ws.Range("N10").ClearContents()
$ws.Range("H16").Value = 1017.5
$ws.Range("I16").Value = 1011.8889
$ws.Range("K16").Value = 1011.8889
$ws.Range("M16").Value = -841.8889
$ws.Range("H55").Value = 34483068
$ws.Range("J55").Value = 489.1875
$ws.Range("L55").Value = 489.1875
$ws.Range("N55").Value = -835.1875
$ws.Range("H68").Value = 3683.5557
$ws.Range("I68").Value = 1692
$ws.Range("K68").Value = 1692
$ws.Range("M68").Value = -943
$ws.Range("H71").Value = 3683.5557
$ws.Range("I71").Value = 1692
$ws.Range("K71").Value = 8460
$ws.Range("M71").Value = -4716
$ws.Range("H130").Value = 59519
$ws.Range("J130").Value = 59519
$ws.Range("L130").Value = 59519
$ws.Range("N130").Value = -69559
$ws.Range("H132").Value = 8338583.5
$ws.Range("I132").Value = 15626903
$ws.Range("K132").Value = 46880709
$ws.Range("M132").Value = -46878179

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 191644.73
$ws.Range("I122").Value = 458831.66
$ws.Range("K122").Value = 1376494.98
$ws.Range("M122").Value = -1374044.98
$ws.Range("H132").Value = 8628432
$ws.Range("I132").Value = 10641392
$ws.Range("J132").Value = 27603.182
$ws.Range("K132").Value = 31924176
$ws.Range("L132").Value = 82809.546
$ws.Range("M132").Value = -31921646
$ws.Range("N132").Value = -87869.546
